$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "JENIS MESIN"
for ($r = 2; $r -le 101; $r++) {
    $ws.Range("I$r").Value = "NAMAMESIN"
}
